# Auto update Excel log
# Appends the latest sensor/alert readings to the ALERTS, Humidity,
# Temperature and Proximity sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# ALERTS sheet: add row 7 (new CRITICAL bathroom alert)
# ---------------------------------------------------------------------
$wsAlerts = $wb.Worksheets.Item("ALERTS")
$alertsRange = $wsAlerts.Range("A7:F7")
$alertsRange.NumberFormat = "@"
$wsAlerts.Cells.Item(7, 1).Value = "2026-02-01"
$wsAlerts.Cells.Item(7, 2).Value = "18:23:17"
$wsAlerts.Cells.Item(7, 3).Value = "18:00"
$wsAlerts.Cells.Item(7, 4).Value = "Bathroom"
$wsAlerts.Cells.Item(7, 5).Value = "CRITICAL"
$wsAlerts.Cells.Item(7, 6).Value = "CRITICAL ALERT: Bathroom occupied, no motion > 60s."
$alertsRange.Style = "Normal"

# ---------------------------------------------------------------------
# Humidity sheet: add rows 55-64
# ---------------------------------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
$humidityRange = $wsHumidity.Range("A55:F64")
$humidityRange.NumberFormat = "@"
$humidityData = @(
  @("2026-02-01","18:23:13","18:00","Bathroom","79.0%","Active"),
  @("2026-02-01","18:23:17","18:00","Bathroom","78.0%","Active"),
  @("2026-02-01","18:23:22","18:00","Bathroom","79.0%","Active"),
  @("2026-02-01","18:23:27","18:00","Bathroom","78.1%","Active"),
  @("2026-02-01","18:23:32","18:00","Bathroom","79.0%","Active"),
  @("2026-02-01","18:23:37","18:00","Bathroom","78.1%","Active"),
  @("2026-02-01","18:23:42","18:00","Bathroom","79.1%","Active"),
  @("2026-02-01","18:23:52","18:00","Bathroom","78.8%","Active"),
  @("2026-02-01","18:24:02","18:00","Bathroom","77.9%","Active"),
  @("2026-02-01","18:24:12","18:00","Bathroom","77.9%","Active")
)
for ($i = 0; $i -lt $humidityData.Count; $i++) {
    $row = $humidityData[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $wsHumidity.Cells.Item(55 + $i, 1 + $j).Value = $row[$j]
    }
}
$humidityRange.Style = "Normal"

# ---------------------------------------------------------------------
# Temperature sheet: add rows 55-64
# ---------------------------------------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
$temperatureRange = $wsTemperature.Range("A55:F64")
$temperatureRange.NumberFormat = "@"
$temperatureData = @(
  @("2026-02-01","18:23:13","18:00","Bathroom","29.4C","Active"),
  @("2026-02-01","18:23:18","18:00","Bathroom","29.3C","Active"),
  @("2026-02-01","18:23:22","18:00","Bathroom","29.4C","Active"),
  @("2026-02-01","18:23:27","18:00","Bathroom","29.3C","Active"),
  @("2026-02-01","18:23:32","18:00","Bathroom","29.4C","Active"),
  @("2026-02-01","18:23:37","18:00","Bathroom","29.4C","Active"),
  @("2026-02-01","18:23:42","18:00","Bathroom","29.4C","Active"),
  @("2026-02-01","18:23:52","18:00","Bathroom","29.4C","Active"),
  @("2026-02-01","18:24:02","18:00","Bathroom","29.4C","Active"),
  @("2026-02-01","18:24:13","18:00","Bathroom","29.4C","Active")
)
for ($i = 0; $i -lt $temperatureData.Count; $i++) {
    $row = $temperatureData[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $wsTemperature.Cells.Item(55 + $i, 1 + $j).Value = $row[$j]
    }
}
$temperatureRange.Style = "Normal"

# ---------------------------------------------------------------------
# Proximity sheet: add row 41 (bathroom door exit)
# ---------------------------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")
$proximityRange = $wsProximity.Range("A41:F41")
$proximityRange.NumberFormat = "@"
$wsProximity.Cells.Item(41, 1).Value = "2026-02-01"
$wsProximity.Cells.Item(41, 2).Value = "18:23:25"
$wsProximity.Cells.Item(41, 3).Value = "18:00"
$wsProximity.Cells.Item(41, 4).Value = "Bathroom Door"
$wsProximity.Cells.Item(41, 5).Value = "EXIT"
$wsProximity.Cells.Item(41, 6).Value = "User EXITED Bathroom"
$proximityRange.Style = "Normal"
